# Update Excel SCD0011 until SCD0016
#
# - Renames the worksheet tab from "SCD0180" to "SCD0011"
# - Updates the TC_ID value in B2 from "DGS-195" to "SCD0011-011"
# - Widens column B so the longer TC_ID text keeps fitting (bestFit-style resize)
# - Moves the active selection from S2 to B3 (also resets the scrolled
#   top-left cell back to the default / A1 corner)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0180 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID value in B2 from "DGS-195" to "SCD0011-011"
$ws.Range("B2").Value = "SCD0011-011"

# Widen column B to fit the new, longer TC_ID text (target ~12.43 chars wide)
$ws.Columns.Item(2).ColumnWidth = 11.67

# Move/reset the active cell selection to B3 and scroll back to the default view
$ws.Range("B3").Select()
